$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty static friction (muStatic) values for the
# "gravel/dry" (row 7) and "sand/dry" (row 8) rows.
$ws.Range("C7").Value = 0.55000000000000004
$ws.Range("C8").Value = 0.35

# Leave the final selection on the last-edited cell, matching the session.
$ws.Range("C8").Select()
